# daily auto push: 2026-02-26 07:14 UTC
#
# A new schedule entry for 2026/02/26 (Thu) got appended to the end of that
# day's block, which lives right before the 2026/12/29 block in the sheet.
# In sheet terms that means: insert a new row at 878 (pushing the old
# 878..919 down to 879..920) and fill it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 878 (and everything below it) down by one to make room.
$ws.Rows(878).Insert()

# Column A holds date strings like "2026/12/29" that must stay literal
# text (not get auto-converted into a date serial number by Excel), so
# mark the cell as Text before writing the value - matches how the rest
# of column A is authored.
$ws.Cells.Item(878, 1).NumberFormat = "@"
$ws.Cells.Item(878, 1).Value = "2026/02/26"
$ws.Cells.Item(878, 2).Value = "木"
$ws.Cells.Item(878, 3).Value = 14
$ws.Cells.Item(878, 4).Value = 201
